# livestock tab1 tidy finished
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table (header row + data rows), matches df-vars-matched.xlsx
# columns: input, chn_block4, asis, variables
$data = @(
    @("总数",     "总数",     $true, "v8_t1_zcqc_zs"),
    @("种牛场",   "种牛场",   $true, "v8_t1_zcqc_znc"),
    @("种奶牛场", "种奶牛场", $true, "v8_t1_zcqc_znnc"),
    @("种肉牛场", "种肉牛场", $true, "v8_t1_zcqc_zrnc"),
    @("种水牛场", "种水牛场", $true, "v8_t1_zcqc_zsnc"),
    @("种牦牛场", "种牦牛场", $true, "v8_t1_zcqc_zhnc"),
    @("种马场",   "种马场",   $true, "v8_t1_zcqc_zmc"),
    @("种猪场",   "种猪场",   $true, "v8_t1_zcqc_zzc"),
    @("种羊场",   "种羊场",   $true, "v8_t1_zcqc_zyc"),
    @("种绵羊场", "种绵羊场", $true, "v8_t1_zcqc_zmyc"),
    @("种细毛羊场", "种细毛羊场", $true, "v8_t1_zcqc_zxmyc"),
    @("种山羊场", "种山羊场", $true, "v8_t1_zcqc_zsyc"),
    @("种绒山羊场", "种绒山羊场", $true, "v8_t1_zcqc_zmsyc"),
    @("种乳牛场", "种乳牛场", $true, "v8_t1_zcqc_zrnc")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}
